$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while forcing text storage
# (prevents Excel from auto-converting numeric-looking strings like
# "209.64" or "0.9990" into real numbers and mangling formatting/precision),
# then restores the cell's style to the workbook default so no stray
# per-cell number-format style is left behind.
function Set-TextValue($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '25.897.00'
Set-TextValue "E2" '  -1.13%  '
Set-TextValue "D3" '1.631.94'
Set-TextValue "E3" '  -2.67%  '
Set-TextValue "E4" '  -0.05%  '
Set-TextValue "D5" '209.64'
Set-TextValue "E5" '  -1.01%  '
Set-TextValue "D6" '0.5199'
Set-TextValue "E6" '  -1.03%  '
Set-TextValue "E7" '  +0.00%  '
Set-TextValue "D8" '0.2567'
Set-TextValue "E8" '  -3.36%  '
Set-TextValue "D9" '0.06224'
Set-TextValue "E9" '  -1.13%  '
Set-TextValue "D10" '20.26'
Set-TextValue "E10" '  -5.29%  '
Set-TextValue "D11" '0.07557'
Set-TextValue "E11" '  -0.02%  '
Set-TextValue "D12" '1.631.49'
Set-TextValue "E12" '  -2.67%  '
Set-TextValue "D13" '4.349'
Set-TextValue "E13" '  -2.37%  '
Set-TextValue "D14" '1.858.84'
Set-TextValue "E14" '  -2.49%  '
Set-TextValue "D15" '0.5414'
Set-TextValue "E15" '  -3.93%  '
Set-TextValue "D16" '0.0₅7908'
Set-TextValue "E16" '  -1.45%  '
Set-TextValue "E17" '  -3.90%  '
Set-TextValue "D18" '25.896.02'
Set-TextValue "E18" '  -1.39%  '
Set-TextValue "E19" '  -0.01%  '
Set-TextValue "D20" '4.602'
Set-TextValue "E20" '  -4.60%  '
Set-TextValue "D21" '184.06'
Set-TextValue "E21" '  -2.14%  '
Set-TextValue "D22" '10.01'
Set-TextValue "E22" '  -4.16%  '
Set-TextValue "D23" '6.060'
Set-TextValue "E23" '  -2.11%  '
Set-TextValue "E24" '  +0.05%  '
Set-TextValue "D25" '145.65'
Set-TextValue "E25" '  -2.57%  '
Set-TextValue "D26" '0.1202'
Set-TextValue "E26" '  -3.86%  '
Set-TextValue "D27" '7.328'
Set-TextValue "E27" '  -3.28%  '
Set-TextValue "D28" '15.46'
Set-TextValue "E28" '  -3.72%  '
Set-TextValue "D29" '1.369'
Set-TextValue "E29" '  +0.75%  '
Set-TextValue "D30" '0.05925'
Set-TextValue "E30" '  -4.39%  '
Set-TextValue "E31" '  -3.28%  '
Set-TextValue "E32" '  -2.81%  '
Set-TextValue "D33" '3.346'
Set-TextValue "E33" '  -4.37%  '
Set-TextValue "D34" '1.602'
Set-TextValue "E34" '  -1.99%  '
Set-TextValue "D35" '0.9691'
Set-TextValue "E35" '  -3.37%  '
Set-TextValue "D36" '2.383'
Set-TextValue "E36" '  -0.99%  '
Set-TextValue "D37" '2.738'
Set-TextValue "E37" '  -0.18%  '
Set-TextValue "D38" '0.5791'
Set-TextValue "E38" '  -4.64%  '
Set-TextValue "D39" '0.01591'
Set-TextValue "E39" '  -1.32%  '
Set-TextValue "E40" '  -0.49%  '
Set-TextValue "D41" '0.8374'
Set-TextValue "E41" '  -3.67%  '
Set-TextValue "D42" '5.643'
Set-TextValue "E42" '  -7.49%  '
Set-TextValue "D43" '1.020.10'
Set-TextValue "E43" '  -6.00%  '
Set-TextValue "D44" '99.52'
Set-TextValue "E44" '  -0.57%  '
Set-TextValue "D45" '1.783.64'
Set-TextValue "E45" '  -2.42%  '
Set-TextValue "D46" '0.0₈108'
Set-TextValue "E46" '  -2.50%  '
Set-TextValue "D47" '0.9990'
Set-TextValue "E47" '  +0.12%  '
Set-TextValue "D48" '54.20'
Set-TextValue "E48" '  -3.65%  '
Set-TextValue "D49" '7.929'
Set-TextValue "E49" '  -0.92%  '
Set-TextValue "D51" '0.4229'
Set-TextValue "E51" '  -0.67%  '

Write-Host "Updated cryptos list"
